$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Greece"
$ws.Cells.Item(2, 2).Value = "April"
$ws.Cells.Item(2, 3).Value = 4.197818542389688
$ws.Cells.Item(2, 4).Value = 2017
$ws.Cells.Item(2, 5).Value = 2.119484382746653
$ws.Cells.Item(2, 6).Value = 132476.5

$ws.Cells.Item(3, 1).Value = "Greece"
$ws.Cells.Item(3, 2).Value = "August"
$ws.Cells.Item(3, 3).Value = 4.351172152498224
$ws.Cells.Item(3, 4).Value = 4223
$ws.Cells.Item(3, 5).Value = 2.178309258820744
$ws.Cells.Item(3, 6).Value = 307807.89

$ws.Cells.Item(4, 1).Value = "Greece"
$ws.Cells.Item(4, 2).Value = "December"
$ws.Cells.Item(4, 3).Value = 4.520833333333333
$ws.Cells.Item(4, 4).Value = 48
$ws.Cells.Item(4, 5).Value = 1.208333333333333
$ws.Cells.Item(4, 6).Value = 2150.8

$ws.Cells.Item(5, 1).Value = "Greece"
$ws.Cells.Item(5, 2).Value = "February"
$ws.Cells.Item(5, 3).Value = 7.565749235474006
$ws.Cells.Item(5, 4).Value = 327
$ws.Cells.Item(5, 5).Value = 1.749235474006116
$ws.Cells.Item(5, 6).Value = 26937

$ws.Cells.Item(6, 1).Value = "Greece"
$ws.Cells.Item(6, 2).Value = "January"
$ws.Cells.Item(6, 3).Value = 4.215568862275449
$ws.Cells.Item(6, 4).Value = 167
$ws.Cells.Item(6, 5).Value = 1.491017964071856
$ws.Cells.Item(6, 6).Value = 7348.5

$ws.Cells.Item(7, 1).Value = "Greece"
$ws.Cells.Item(7, 2).Value = "June"
$ws.Cells.Item(7, 3).Value = 4.164274924471299
$ws.Cells.Item(7, 4).Value = 2648
$ws.Cells.Item(7, 5).Value = 1.95392749244713
$ws.Cells.Item(7, 6).Value = 182864

$ws.Cells.Item(8, 1).Value = "Greece"
$ws.Cells.Item(8, 2).Value = "March"
$ws.Cells.Item(8, 3).Value = 5.923076923076923
$ws.Cells.Item(8, 4).Value = 572
$ws.Cells.Item(8, 5).Value = 1.694055944055944
$ws.Cells.Item(8, 6).Value = 34128

$ws.Cells.Item(9, 1).Value = "Greece"
$ws.Cells.Item(9, 2).Value = "May"
$ws.Cells.Item(9, 3).Value = 4.27736890524379
$ws.Cells.Item(9, 4).Value = 2174
$ws.Cells.Item(9, 5).Value = 2.10395584176633
$ws.Cells.Item(9, 6).Value = 150305.46

$ws.Cells.Item(10, 1).Value = "Greece"
$ws.Cells.Item(10, 2).Value = "November"
$ws.Cells.Item(10, 3).Value = 4.117021276595745
$ws.Cells.Item(10, 4).Value = 282
$ws.Cells.Item(10, 5).Value = 1.439716312056738
$ws.Cells.Item(10, 6).Value = 12941.32

$ws.Cells.Item(11, 1).Value = "Greece"
$ws.Cells.Item(11, 2).Value = "October"
$ws.Cells.Item(11, 3).Value = 4.22294776119403
$ws.Cells.Item(11, 4).Value = 1072
$ws.Cells.Item(11, 5).Value = 1.800373134328358
$ws.Cells.Item(11, 6).Value = 72117.32000000001

$ws.Cells.Item(12, 1).Value = "Greece"
$ws.Cells.Item(12, 2).Value = "September"
$ws.Cells.Item(12, 3).Value = 4.158210947930574
$ws.Cells.Item(12, 4).Value = 1498
$ws.Cells.Item(12, 5).Value = 2.068758344459279
$ws.Cells.Item(12, 6).Value = 102478.57

$ws.Cells.Item(13, 1).Value = "Italy"
$ws.Cells.Item(13, 2).Value = "April"
$ws.Cells.Item(13, 3).Value = 4.205128205128205
$ws.Cells.Item(13, 4).Value = 312
$ws.Cells.Item(13, 5).Value = 1.375
$ws.Cells.Item(13, 6).Value = 28415

$ws.Cells.Item(14, 1).Value = "Italy"
$ws.Cells.Item(14, 2).Value = "August"
$ws.Cells.Item(14, 3).Value = 3
$ws.Cells.Item(14, 4).Value = 40
$ws.Cells.Item(14, 5).Value = 1.375
$ws.Cells.Item(14, 6).Value = 2245

$ws.Cells.Item(15, 1).Value = "Italy"
$ws.Cells.Item(15, 2).Value = "February"
$ws.Cells.Item(15, 3).Value = 3.461538461538462
$ws.Cells.Item(15, 4).Value = 104
$ws.Cells.Item(15, 5).Value = 1.375
$ws.Cells.Item(15, 6).Value = 7229

$ws.Cells.Item(16, 1).Value = "Italy"
$ws.Cells.Item(16, 2).Value = "January"
$ws.Cells.Item(16, 3).Value = 4
$ws.Cells.Item(16, 4).Value = 120
$ws.Cells.Item(16, 5).Value = 1.375
$ws.Cells.Item(16, 6).Value = 10215

$ws.Cells.Item(17, 1).Value = "Italy"
$ws.Cells.Item(17, 2).Value = "June"
$ws.Cells.Item(17, 3).Value = 3.777777777777778
$ws.Cells.Item(17, 4).Value = 72
$ws.Cells.Item(17, 5).Value = 1.375
$ws.Cells.Item(17, 6).Value = 5665

$ws.Cells.Item(18, 1).Value = "Italy"
$ws.Cells.Item(18, 2).Value = "March"
$ws.Cells.Item(18, 3).Value = 6
$ws.Cells.Item(18, 4).Value = 88
$ws.Cells.Item(18, 5).Value = 1.375
$ws.Cells.Item(18, 6).Value = 12595

$ws.Cells.Item(19, 1).Value = "Italy"
$ws.Cells.Item(19, 2).Value = "May"
$ws.Cells.Item(19, 3).Value = 5.272727272727272
$ws.Cells.Item(19, 4).Value = 88
$ws.Cells.Item(19, 5).Value = 1.375
$ws.Cells.Item(19, 6).Value = 10739

$ws.Cells.Item(20, 1).Value = "Italy"
$ws.Cells.Item(20, 2).Value = "November"
$ws.Cells.Item(20, 3).Value = 4.557692307692307
$ws.Cells.Item(20, 4).Value = 52
$ws.Cells.Item(20, 5).Value = 1.288461538461539
$ws.Cells.Item(20, 6).Value = 4614

$ws.Cells.Item(21, 1).Value = "Italy"
$ws.Cells.Item(21, 2).Value = "October"
$ws.Cells.Item(21, 3).Value = 3.844660194174757
$ws.Cells.Item(21, 4).Value = 206
$ws.Cells.Item(21, 5).Value = 1.276699029126214
$ws.Cells.Item(21, 6).Value = 13953

$ws.Cells.Item(22, 1).Value = "Italy"
$ws.Cells.Item(22, 2).Value = "September"
$ws.Cells.Item(22, 3).Value = 3.75
$ws.Cells.Item(22, 4).Value = 64
$ws.Cells.Item(22, 5).Value = 1.375
$ws.Cells.Item(22, 6).Value = 4161.12

$ws.Cells.Item(23, 1).Value = "Portugal"
$ws.Cells.Item(23, 2).Value = "April"
$ws.Cells.Item(23, 3).Value = 4.089038385437277
$ws.Cells.Item(23, 4).Value = 2527
$ws.Cells.Item(23, 5).Value = 1.6398891966759
$ws.Cells.Item(23, 6).Value = 173967.5

$ws.Cells.Item(24, 1).Value = "Portugal"
$ws.Cells.Item(24, 2).Value = "August"
$ws.Cells.Item(24, 3).Value = 4.164769647696477
$ws.Cells.Item(24, 4).Value = 1845
$ws.Cells.Item(24, 5).Value = 1.38319783197832
$ws.Cells.Item(24, 6).Value = 119299.68

$ws.Cells.Item(25, 1).Value = "Portugal"
$ws.Cells.Item(25, 2).Value = "December"
$ws.Cells.Item(25, 3).Value = 4.112244897959184
$ws.Cells.Item(25, 4).Value = 196
$ws.Cells.Item(25, 5).Value = 1.173469387755102
$ws.Cells.Item(25, 6).Value = 10991.2

$ws.Cells.Item(26, 1).Value = "Portugal"
$ws.Cells.Item(26, 2).Value = "February"
$ws.Cells.Item(26, 3).Value = 3.720812182741117
$ws.Cells.Item(26, 4).Value = 197
$ws.Cells.Item(26, 5).Value = 1.467005076142132
$ws.Cells.Item(26, 6).Value = 10560.04

$ws.Cells.Item(27, 1).Value = "Portugal"
$ws.Cells.Item(27, 2).Value = "January"
$ws.Cells.Item(27, 3).Value = 3.961832061068702
$ws.Cells.Item(27, 4).Value = 131
$ws.Cells.Item(27, 5).Value = 1.595419847328244
$ws.Cells.Item(27, 6).Value = 7772

$ws.Cells.Item(28, 1).Value = "Portugal"
$ws.Cells.Item(28, 2).Value = "June"
$ws.Cells.Item(28, 3).Value = 4.231991525423729
$ws.Cells.Item(28, 4).Value = 1888
$ws.Cells.Item(28, 5).Value = 1.621822033898305
$ws.Cells.Item(28, 6).Value = 136740.7

$ws.Cells.Item(29, 1).Value = "Portugal"
$ws.Cells.Item(29, 2).Value = "March"
$ws.Cells.Item(29, 3).Value = 3.949392712550607
$ws.Cells.Item(29, 4).Value = 494
$ws.Cells.Item(29, 5).Value = 1.629554655870445
$ws.Cells.Item(29, 6).Value = 33100.7

$ws.Cells.Item(30, 1).Value = "Portugal"
$ws.Cells.Item(30, 2).Value = "May"
$ws.Cells.Item(30, 3).Value = 4.262116316639742
$ws.Cells.Item(30, 4).Value = 2476
$ws.Cells.Item(30, 5).Value = 1.644184168012924
$ws.Cells.Item(30, 6).Value = 179698

$ws.Cells.Item(31, 1).Value = "Portugal"
$ws.Cells.Item(31, 2).Value = "November"
$ws.Cells.Item(31, 3).Value = 4.085106382978723
$ws.Cells.Item(31, 4).Value = 329
$ws.Cells.Item(31, 5).Value = 1.288753799392097
$ws.Cells.Item(31, 6).Value = 20496.6

$ws.Cells.Item(32, 1).Value = "Portugal"
$ws.Cells.Item(32, 2).Value = "October"
$ws.Cells.Item(32, 3).Value = 3.954545454545455
$ws.Cells.Item(32, 4).Value = 1100
$ws.Cells.Item(32, 5).Value = 1.493636363636364
$ws.Cells.Item(32, 6).Value = 70794.60000000001

$ws.Cells.Item(33, 1).Value = "Portugal"
$ws.Cells.Item(33, 2).Value = "September"
$ws.Cells.Item(33, 3).Value = 3.888774459320288
$ws.Cells.Item(33, 4).Value = 971
$ws.Cells.Item(33, 5).Value = 1.440782698249228
$ws.Cells.Item(33, 6).Value = 58716

$ws.Cells.Item(34, 1).Value = "Spain"
$ws.Cells.Item(34, 2).Value = "April"
$ws.Cells.Item(34, 3).Value = 4.551724137931035
$ws.Cells.Item(34, 4).Value = 87
$ws.Cells.Item(34, 5).Value = 1.172413793103448
$ws.Cells.Item(34, 6).Value = 7231.46

$ws.Cells.Item(35, 1).Value = "Spain"
$ws.Cells.Item(35, 2).Value = "August"
$ws.Cells.Item(35, 3).Value = 4.133928571428571
$ws.Cells.Item(35, 4).Value = 224
$ws.Cells.Item(35, 5).Value = 1.209821428571429
$ws.Cells.Item(35, 6).Value = 18747.7

$ws.Cells.Item(36, 1).Value = "Spain"
$ws.Cells.Item(36, 2).Value = "December"
$ws.Cells.Item(36, 3).Value = 3.777777777777778
$ws.Cells.Item(36, 4).Value = 45
$ws.Cells.Item(36, 5).Value = 1.155555555555555
$ws.Cells.Item(36, 6).Value = 2890

$ws.Cells.Item(37, 1).Value = "Spain"
$ws.Cells.Item(37, 2).Value = "June"
$ws.Cells.Item(37, 3).Value = 3.567901234567901
$ws.Cells.Item(37, 4).Value = 81
$ws.Cells.Item(37, 5).Value = 1.185185185185185
$ws.Cells.Item(37, 6).Value = 5520.82

$ws.Cells.Item(38, 1).Value = "Spain"
$ws.Cells.Item(38, 2).Value = "May"
$ws.Cells.Item(38, 3).Value = 3.777777777777778
$ws.Cells.Item(38, 4).Value = 81
$ws.Cells.Item(38, 5).Value = 1.209876543209877
$ws.Cells.Item(38, 6).Value = 5793.78

$ws.Cells.Item(39, 1).Value = "Spain"
$ws.Cells.Item(39, 2).Value = "November"
$ws.Cells.Item(39, 3).Value = 4.212355212355212
$ws.Cells.Item(39, 4).Value = 259
$ws.Cells.Item(39, 5).Value = 1.173745173745174
$ws.Cells.Item(39, 6).Value = 19154.92

$ws.Cells.Item(40, 1).Value = "Spain"
$ws.Cells.Item(40, 2).Value = "October"
$ws.Cells.Item(40, 3).Value = 3.787446504992867
$ws.Cells.Item(40, 4).Value = 701
$ws.Cells.Item(40, 5).Value = 1.156918687589158
$ws.Cells.Item(40, 6).Value = 47079.84

$ws.Cells.Item(41, 1).Value = "Spain"
$ws.Cells.Item(41, 2).Value = "September"
$ws.Cells.Item(41, 3).Value = 3.989690721649485
$ws.Cells.Item(41, 4).Value = 194
$ws.Cells.Item(41, 5).Value = 1.164948453608248
$ws.Cells.Item(41, 6).Value = 13759.46

$ws.Cells.Item(42, 1).Value = "Turkey"
$ws.Cells.Item(42, 2).Value = "October"
$ws.Cells.Item(42, 3).Value = 2.311111111111111
$ws.Cells.Item(42, 4).Value = 45
$ws.Cells.Item(42, 5).Value = 1
$ws.Cells.Item(42, 6).Value = 2038

$ws.Cells.Item(43, 1).Value = "Turkey"
$ws.Cells.Item(43, 2).Value = "September"
$ws.Cells.Item(43, 3).Value = 2.64
$ws.Cells.Item(43, 4).Value = 50
$ws.Cells.Item(43, 5).Value = 1
$ws.Cells.Item(43, 6).Value = 2710

$ws.Range("A44:F47").ClearContents()
